$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.960.33'
$ws.Range('E2').Value = '  +5.24%  '
$ws.Range('D3').Value = '2.292.07'
$ws.Range('E3').Value = '  +2.72%  '
$ws.Range('E4').Value = '  +0.14%  '
$__style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '232.19'
$ws.Range('D5').Style = $__style
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('E6').Value = '  +1.65%  '
$__style = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '63.84'
$ws.Range('D7').Style = $__style
$ws.Range('E7').Value = '  +5.29%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  +5.06%  '
$__style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0952'
$ws.Range('D10').Style = $__style
$ws.Range('E10').Value = '  +4.99%  '
$__style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '57.64'
$ws.Range('D11').Style = $__style
$ws.Range('E11').Value = '  -0.82%  '
$__style = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '26.33'
$ws.Range('D12').Style = $__style
$ws.Range('E12').Value = '  +15.31%  '
$ws.Range('E13').Value = '  +0.33%  '
$ws.Range('D14').Value = '2.633.12'
$ws.Range('E14').Value = '  +2.75%  '
$__style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.88'
$ws.Range('D15').Style = $__style
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('E16').Value = '  +6.12%  '
$__style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.817'
$ws.Range('D17').Style = $__style
$ws.Range('E17').Value = '  +1.58%  '
$ws.Range('D18').Value = '2.287.99'
$ws.Range('E18').Value = '  +2.32%  '
$ws.Range('D19').Value = '43.850.16'
$ws.Range('E19').Value = '  +5.05%  '
$ws.Range('D20').Value = '0.0₃0949'
$ws.Range('E20').Value = '  +4.70%  '
$__style = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '73.34'
$ws.Range('D21').Style = $__style
$ws.Range('E21').Value = '  +1.24%  '
$__style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.20'
$ws.Range('D22').Style = $__style
$ws.Range('E22').Value = '  +1.67%  '
$__style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '251.95'
$ws.Range('D23').Style = $__style
$ws.Range('E23').Value = '  +1.43%  '
$__style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.60'
$ws.Range('D24').Style = $__style
$ws.Range('E24').Value = '  +8.97%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('E26').Value = '  +0.96%  '
$__style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.91'
$ws.Range('D27').Style = $__style
$ws.Range('E27').Value = '  +1.54%  '
$__style = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '171.67'
$ws.Range('D28').Style = $__style
$ws.Range('E28').Value = '  +1.16%  '
$__style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.140'
$ws.Range('D29').Style = $__style
$ws.Range('E29').Value = '  -2.15%  '
$ws.Range('E30').Value = '  +3.25%  '
$__style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.44'
$ws.Range('D31').Style = $__style
$ws.Range('E31').Value = '  +2.37%  '
$ws.Range('E32').Value = '  +4.02%  '
$ws.Range('E33').Value = '  +0.33%  '
$__style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0704'
$ws.Range('D34').Style = $__style
$ws.Range('E34').Value = '  +7.77%  '
$__style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.15'
$ws.Range('D35').Style = $__style
$ws.Range('E35').Value = '  +1.59%  '
$ws.Range('E36').Value = '  +0.80%  '
$ws.Range('B37').Value = 'THORChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$__style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.62'
$ws.Range('D37').Style = $__style
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$__style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.71'
$ws.Range('D38').Style = $__style
$ws.Range('E38').Value = '  +1.71%  '
$__style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.37'
$ws.Range('D39').Style = $__style
$ws.Range('E39').Value = '  -0.87%  '
$__style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0249'
$ws.Range('D40').Style = $__style
$ws.Range('E40').Value = '  +3.47%  '
$ws.Range('E41').Value = '  -0.06%  '
$__style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.11'
$ws.Range('D42').Style = $__style
$ws.Range('E42').Value = '  +28.03%  '
$__style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.58'
$ws.Range('D43').Style = $__style
$ws.Range('E43').Value = '  +0.09%  '
$__style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.64'
$ws.Range('D44').Style = $__style
$ws.Range('E44').Value = '  +3.77%  '
$ws.Range('E45').Value = '  -7.69%  '
$ws.Range('E46').Value = '  -0.23%  '
$__style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0966'
$ws.Range('D47').Style = $__style
$ws.Range('E47').Value = '  +0.72%  '
$__style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '98.45'
$ws.Range('D48').Style = $__style
$ws.Range('E48').Value = '  -0.26%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '1.490.48'
$ws.Range('E49').Value = '  +1.46%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$__style = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '17.02'
$ws.Range('D50').Style = $__style
$ws.Range('E50').Value = '  +2.51%  '
$ws.Range('E51').Value = '  +3.33%  '
